$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 204; existing rows 204-212 shift
# down to 205-213 (data unchanged), and this new row 204 receives the
# latest weekly price observation.
$ws.Rows.Item(204).Insert()

$ws.Cells.Item(204, 1).Value = 7
$ws.Cells.Item(204, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(204, 3).Value = "Ñuble"
$ws.Cells.Item(204, 4).Value = 45075
$ws.Cells.Item(204, 5).Value = 16
$ws.Cells.Item(204, 6).Value = 100112021
$ws.Cells.Item(204, 7).Value = "Ají"
$ws.Cells.Item(204, 8).Value = "Cacho cabra rojo"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 30
$ws.Cells.Item(204, 11).Value = 17000
$ws.Cells.Item(204, 12).Value = 17000
$ws.Cells.Item(204, 13).Value = 17000
$ws.Cells.Item(204, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(204, 15).Value = "Región del Maule"
$ws.Cells.Item(204, 16).Value = 680
$ws.Cells.Item(204, 17).Value = 25
$ws.Cells.Item(204, 18).Value = "Hortaliza"
